$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Remove the old aggregate sheet first so the id counter frees up "2"
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Delete()

# Duplicate "2021-Q4" right after itself -> becomes the new quarter sheet.
# This preserves the exact header/column-A styling (bold, centered, thin border)
# that the original sheet already carries.
$sheet1.Copy($null, $sheet1)
$ws2022 = $wb.Worksheets.Item(2)
$ws2022.Name = "2022-Q1"


# The copied sheet has 19 data rows (rows 2-20); the 2022-Q1 snapshot only
# has 17 holdings, so drop the trailing two rows before writing fresh data.
$ws2022.Range("A19:A20").EntireRow.Delete()

# Overwrite the 17 data rows with the 2022-Q1 fund holdings.
# Numeric-looking text columns (fund code / scale / position figures) are
# entered with a leading apostrophe so they stay text (keeps leading zeros
# like "007490" and trailing zeros like "21.90" intact).
$ws2022.Range("A2").Value = 0
$ws2022.Range("B2").Value = "'007490"
$ws2022.Range("C2").Value = "南方信息创新混合A"
$ws2022.Range("D2").Value = "'21.90"
$ws2022.Range("E2").Value = "'91.22"
$ws2022.Range("F2").Value = "'4.11"
$ws2022.Range("G2").Value = "'0.9001"
$ws2022.Range("H2").Value = 8
$ws2022.Range("A3").Value = 1
$ws2022.Range("B3").Value = "'001404"
$ws2022.Range("C3").Value = "招商移动互联网产业股票"
$ws2022.Range("D3").Value = "'13.45"
$ws2022.Range("E3").Value = "'90.96"
$ws2022.Range("F3").Value = "'5.18"
$ws2022.Range("G3").Value = "'0.6967"
$ws2022.Range("H3").Value = 6
$ws2022.Range("A4").Value = 2
$ws2022.Range("B4").Value = "'001042"
$ws2022.Range("C4").Value = "华夏领先股票"
$ws2022.Range("D4").Value = "'13.72"
$ws2022.Range("E4").Value = "'93.46"
$ws2022.Range("F4").Value = "'2.51"
$ws2022.Range("G4").Value = "'0.3444"
$ws2022.Range("H4").Value = 7
$ws2022.Range("A5").Value = 3
$ws2022.Range("B5").Value = "'008655"
$ws2022.Range("C5").Value = "招商科技创新混合A"
$ws2022.Range("D5").Value = "'4.72"
$ws2022.Range("E5").Value = "'90.71"
$ws2022.Range("F5").Value = "'5.25"
$ws2022.Range("G5").Value = "'0.2478"
$ws2022.Range("H5").Value = 3
$ws2022.Range("A6").Value = 4
$ws2022.Range("B6").Value = "'012556"
$ws2022.Range("C6").Value = "长盛景气优选混合"
$ws2022.Range("D6").Value = "'10.87"
$ws2022.Range("E6").Value = "'47.88"
$ws2022.Range("F6").Value = "'1.31"
$ws2022.Range("G6").Value = "'0.1424"
$ws2022.Range("H6").Value = 7
$ws2022.Range("A7").Value = 5
$ws2022.Range("B7").Value = "'007491"
$ws2022.Range("C7").Value = "南方信息创新混合C"
$ws2022.Range("D7").Value = "'2.59"
$ws2022.Range("E7").Value = "'91.22"
$ws2022.Range("F7").Value = "'4.11"
$ws2022.Range("G7").Value = "'0.1064"
$ws2022.Range("H7").Value = 8
$ws2022.Range("A8").Value = 6
$ws2022.Range("B8").Value = "'004314"
$ws2022.Range("C8").Value = "前海开源沪港深新硬件主题灵活配置混合A"
$ws2022.Range("D8").Value = "'1.67"
$ws2022.Range("E8").Value = "'90.05"
$ws2022.Range("F8").Value = "'5.15"
$ws2022.Range("G8").Value = "'0.0860"
$ws2022.Range("H8").Value = 6
$ws2022.Range("A9").Value = 7
$ws2022.Range("B9").Value = "'008656"
$ws2022.Range("C9").Value = "招商科技创新混合C"
$ws2022.Range("D9").Value = "'1.57"
$ws2022.Range("E9").Value = "'90.71"
$ws2022.Range("F9").Value = "'5.25"
$ws2022.Range("G9").Value = "'0.0824"
$ws2022.Range("H9").Value = 3
$ws2022.Range("A10").Value = 8
$ws2022.Range("B10").Value = "'012200"
$ws2022.Range("C10").Value = "新华鑫科技3个月滚动持有灵活配置混合型证券投资基金A"
$ws2022.Range("D10").Value = "'2.04"
$ws2022.Range("E10").Value = "'77.02"
$ws2022.Range("F10").Value = "'3.71"
$ws2022.Range("G10").Value = "'0.0757"
$ws2022.Range("H10").Value = 3
$ws2022.Range("A11").Value = 9
$ws2022.Range("B11").Value = "'013339"
$ws2022.Range("C11").Value = "创金合信芯片产业股票A"
$ws2022.Range("D11").Value = "'1.50"
$ws2022.Range("E11").Value = "'93.43"
$ws2022.Range("F11").Value = "'4.64"
$ws2022.Range("G11").Value = "'0.0696"
$ws2022.Range("H11").Value = 9
$ws2022.Range("A12").Value = 10
$ws2022.Range("B12").Value = "'004315"
$ws2022.Range("C12").Value = "前海开源沪港深新硬件主题灵活配置混合C"
$ws2022.Range("D12").Value = "'1.00"
$ws2022.Range("E12").Value = "'90.05"
$ws2022.Range("F12").Value = "'5.15"
$ws2022.Range("G12").Value = "'0.0515"
$ws2022.Range("H12").Value = 6
$ws2022.Range("A13").Value = 11
$ws2022.Range("B13").Value = "'004044"
$ws2022.Range("C13").Value = "金鹰转型动力灵活配置混合"
$ws2022.Range("D13").Value = "'0.72"
$ws2022.Range("E13").Value = "'93.34"
$ws2022.Range("F13").Value = "'5.18"
$ws2022.Range("G13").Value = "'0.0373"
$ws2022.Range("H13").Value = 7
$ws2022.Range("A14").Value = 12
$ws2022.Range("B14").Value = "'013340"
$ws2022.Range("C14").Value = "创金合信芯片产业股票C"
$ws2022.Range("D14").Value = "'0.61"
$ws2022.Range("E14").Value = "'93.43"
$ws2022.Range("F14").Value = "'4.64"
$ws2022.Range("G14").Value = "'0.0283"
$ws2022.Range("H14").Value = 9
$ws2022.Range("A15").Value = 13
$ws2022.Range("B15").Value = "'012201"
$ws2022.Range("C15").Value = "新华鑫科技3个月滚动持有灵活配置混合型证券投资基金C"
$ws2022.Range("D15").Value = "'0.52"
$ws2022.Range("E15").Value = "'77.02"
$ws2022.Range("F15").Value = "'3.71"
$ws2022.Range("G15").Value = "'0.0193"
$ws2022.Range("H15").Value = 3
$ws2022.Range("A16").Value = 14
$ws2022.Range("B16").Value = "'001574"
$ws2022.Range("C16").Value = "中海混改红利主题精选灵活配置混合"
$ws2022.Range("D16").Value = "'0.30"
$ws2022.Range("E16").Value = "'89.77"
$ws2022.Range("F16").Value = "'3.69"
$ws2022.Range("G16").Value = "'0.0111"
$ws2022.Range("H16").Value = 9
$ws2022.Range("A17").Value = 15
$ws2022.Range("B17").Value = "'002303"
$ws2022.Range("C17").Value = "金鹰智慧生活灵活配置混合"
$ws2022.Range("D17").Value = "'0.11"
$ws2022.Range("E17").Value = "'89.88"
$ws2022.Range("F17").Value = "'6.76"
$ws2022.Range("G17").Value = "'0.0074"
$ws2022.Range("H17").Value = 4
$ws2022.Range("A18").Value = 16
$ws2022.Range("B18").Value = "'161721"
$ws2022.Range("C18").Value = "招商沪深300地产等权重指数"
$ws2022.Range("D18").Value = "'9.97"
$ws2022.Range("E18").Value = "'94.51"
$ws2022.Range("F18").Value = "'0.03"
$ws2022.Range("G18").Value = "'0.0030"
$ws2022.Range("H18").Value = 10

# Build the new "总计" (totals) sheet by duplicating the just-created
# 2022-Q1 sheet (so it inherits the same s="2" bold/border styling), then
# trim it down to the small 4-column summary table shape.
$ws2022.Copy($null, $ws2022)
$wsTotal = $wb.Worksheets.Item(3)
$wsTotal.Name = "总计"

$wsTotal.Range("A4:A18").EntireRow.Delete()
$wsTotal.Range("E1:H1").EntireColumn.Delete()

$wsTotal.Range("B1").Value = "日期"
$wsTotal.Range("C1").Value = "持有数量(只)"
$wsTotal.Range("D1").Value = "持有市值(亿元)"

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 17
$wsTotal.Range("D2").Value = 2.91

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 19
$wsTotal.Range("D3").Value = 8.57

Write-Host "Sheets: $($wb.Worksheets | ForEach-Object { $_.Name })"
